$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '44.193.64'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.227.61'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -1.82%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '298.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '90.67'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.33%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.559'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.496'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.04%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.28'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.41%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0781'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.01'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.565.02'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.223.13'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '13.41'
$ws.Range('D16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.779'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -6.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '44.026.81'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.29'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.70%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.59%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.99'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.42%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '64.26'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.99'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('E24').Value = '  -4.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.85'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -7.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.27'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '38.87'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.40'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.25%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.25'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '151.40'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.51'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -8.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0763'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.15%  '
$ws.Range('E34').Value = '  -6.18%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.117'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.104'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -5.77%  '
$ws.Range('E37').Value = '  -7.77%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.69'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.97%  '
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.17'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.66%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.60'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '13.52'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -9.34%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.797.79'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.83'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +10.94%  '
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '67.98'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.74%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '94.79'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.28%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '73.31'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -7.28%  '
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.62'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.85%  '
